$d = $word.ActiveDocument

# --- Change 1: Title line - merge "18_" + "PA" + "S ON DHCP Configuration" into "18_PAS ON DHCP Configuration"
#     (also drops the stray w:color="000000" runs because Find/Replace uses the first matched run's formatting)
$d.Content.Find.Execute(
    "DV162_18_PAS ON DHCP Configuration",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "DV162_18_PAS ON DHCP Configuration",
    2)

# --- Change 2: "Possible Answer " + "Sheet" -> "Possible Answer Sheet"
$d.Content.Find.Execute(
    "Possible Answer Sheet",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Possible Answer Sheet",
    2)

# --- Change 3: dynamic assignment definition - merge 3 runs (drop bold on "for a temporary period")
$d.Content.Find.Execute(
    "A dynamic assignment refers to the process of automatically allocating IP addresses to devices on a network for a temporary period. This is  the core function of DHCP compared to static IP assignment.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "A dynamic assignment refers to the process of automatically allocating IP addresses to devices on a network for a temporary period. This is  the core function of DHCP compared to static IP assignment.",
    2)

# --- Change 4: automatic assignment definition - merge 3 runs (drop bold on "for a temporary period")
$d.Content.Find.Execute(
    "the process of  automatically  distributing IP addresses to devices on a network for a temporary period.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "the process of  automatically  distributing IP addresses to devices on a network for a temporary period.",
    2)

# --- Change 5: T2 timer definition - merge 3 runs (drop bold on "DHCP server")
$d.Content.Find.Execute(
    "It defines the timeframe after which a device attempts to obtain a new IP address from DHCP server if it fails to renew its lease with the original server that granted it the lease initially. T2 T",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "It defines the timeframe after which a device attempts to obtain a new IP address from DHCP server if it fails to renew its lease with the original server that granted it the lease initially. T2 T",
    2)

# --- Change 6: merge the final 4 paragraphs (Q7 answer) into a single paragraph / single run
$p1 = $d.Paragraphs.Item(51)
$p2 = $d.Paragraphs.Item(52)
$p3 = $d.Paragraphs.Item(53)

# Delete the paragraph marks joining paragraph 51->52->53->54 into one paragraph.
# (delete from the end backwards so the ranges/indexes found above stay valid)
$mark3 = $d.Range($p3.Range.End - 1, $p3.Range.End)
$mark3.Delete()
$mark2 = $d.Range($p2.Range.End - 1, $p2.Range.End)
$mark2.Delete()
$mark1 = $d.Range($p1.Range.End - 1, $p1.Range.End)
$mark1.Delete()

# Re-join the text with the spaces that used to separate the paragraphs, and collapse to one run.
$d.Content.Find.Execute(
    "A. Generally, the recommended time to lease an IP address is 48 hours to renew the IP addressonce a day. After applying the specified parameters, clients will receive an IP address for 1minute, after which they will send a request to the DHCP server for a new IP address every 30seconds.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "A. Generally, the recommended time to lease an IP address is 48 hours to renew the IP address once a day. After applying the specified parameters, clients will receive an IP address for 1 minute, after which they will send a request to the DHCP server for a new IP address every 30 seconds.",
    2)
